$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(117, 0, $false, 0.953125),
    @(118, 1, $false, 0.96875),
    @(119, 2, $false, 0.953125),
    @(120, 3, $false, 0.953125),
    @(121, 4, $false, 0.921875),
    @(122, 5, $false, 0.890625),
    @(123, 6, $false, 0.875),
    @(124, 7, $false, 0.765625),
    @(125, 8, $false, 0.9375),
    @(126, 9, $false, 0.84375),
    @(127, 10, $false, 0.859375),
    @(128, 11, $false, 0.8125),
    @(129, 12, $false, 0.828125),
    @(130, 13, $false, 0.796875),
    @(131, 14, $false, 0.84375),
    @(132, 15, $false, 0.796875),
    @(133, 16, $false, 0.765625),
    @(134, 17, $false, 0.828125),
    @(135, 18, $false, 0.796875),
    @(136, 19, $false, 0.734375),
    @(137, 20, $false, 0.765625),
    @(138, 21, $false, 0.75),
    @(139, 22, $false, 0.75),
    @(140, 23, $false, 0.765625),
    @(141, 24, $false, 0.75),
    @(142, 25, $false, 0.75),
    @(143, 26, $false, 0.75),
    @(144, 27, $false, 0.75),
    @(145, 28, $false, 0.75),
    @(146, 29, $false, 0.75),
    @(147, 30, $false, 0.734375),
    @(148, 31, $false, 0.734375),
    @(149, 32, $false, 0.734375),
    @(150, 33, $false, 0.734375),
    @(151, 34, $false, 0.75),
    @(152, 35, $false, 0.75),
    @(153, 36, $false, 0.75),
    @(154, 37, $false, 0.765625),
    @(155, 38, $false, 0.78125),
    @(156, 39, $false, 0.78125),
    @(157, 40, $false, 0.78125),
    @(158, 41, $false, 0.78125),
    @(159, 42, $false, 0.796875),
    @(160, 43, $false, 0.796875),
    @(161, 44, $false, 0.796875),
    @(162, 45, $false, 0.796875),
    @(163, 46, $false, 0.796875),
    @(164, 47, $false, 0.796875),
    @(165, 48, $false, 0.78125),
    @(166, 49, $false, 0.78125),
    @(167, 50, $false, 0.78125),
    @(168, 51, $false, 0.78125),
    @(169, 52, $false, 0.78125),
    @(170, 53, $false, 0.78125),
    @(171, 54, $false, 0.78125),
    @(172, 55, $false, 0.78125),
    @(173, 56, $false, 0.78125),
    @(174, 57, $false, 0.78125),
    @(175, 58, $false, 0.765625),
    @(176, 59, $false, 0.765625),
    @(177, 60, $false, 0.75),
    @(178, 61, $false, 0.765625),
    @(179, 62, $false, 0.75),
    @(180, 63, $false, 0.75),
    @(181, 64, $false, 0.75),
    @(182, 65, $false, 0.75),
    @(183, 66, $false, 0.75),
    @(184, 67, $false, 0.75),
    @(185, 68, $false, 0.75),
    @(186, 69, $false, 0.75),
    @(187, 70, $false, 0.75),
    @(188, 71, $false, 0.75),
    @(189, 72, $false, 0.75),
    @(190, 73, $false, 0.75),
    @(191, 74, $false, 0.75),
    @(192, 75, $false, 0.75),
    @(193, 76, $false, 0.75),
    @(194, 77, $false, 0.75),
    @(195, 78, $false, 0.75),
    @(196, 79, $false, 0.75),
    @(197, 80, $false, 0.75),
    @(198, 81, $false, 0.75),
    @(199, 82, $false, 0.75),
    @(200, 83, $false, 0.75),
    @(201, 84, $false, 0.75),
    @(202, 85, $false, 0.75),
    @(203, 86, $false, 0.75),
    @(204, 87, $false, 0.75),
    @(205, 88, $false, 0.75),
    @(206, 89, $false, 0.75),
    @(207, 90, $false, 0.75),
    @(208, 91, $false, 0.75),
    @(209, 92, $false, 0.75),
    @(210, 93, $false, 0.75),
    @(211, 94, $false, 0.75),
    @(212, 95, $false, 0.75),
    @(213, 96, $false, 0.75),
    @(214, 97, $false, 0.75),
    @(215, 98, $false, 0.75),
    @(216, 99, $false, 0.75),
    @(217, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.75),
    @(218, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.796875),
    @(219, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.890625),
    @(220, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.828125),
    @(221, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.796875),
    @(222, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.734375),
    @(223, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.765625),
    @(224, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.75),
    @(225, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.828125),
    @(226, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.75),
    @(227, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.859375),
    @(228, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.921875),
    @(229, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.875),
    @(230, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.875),
    @(231, '<__main__.DisplayOutputs object at 0x7fcc30090610>', $true, 0.7941176470588235)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $aVal = $r[1]
    $isStr = $r[2]
    $bVal = $r[3]

    if ($isStr) {
        $ws.Cells.Item($rowNum, 1).Value = [string]$aVal
    } else {
        $ws.Cells.Item($rowNum, 1).Value = $aVal
    }
    $ws.Cells.Item($rowNum, 2).Value = $bVal
}
